$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.048.73"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.833.14"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'242.75"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "'0.6184"
$ws.Range("E6").Value = "  -2.28%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.07461"
$ws.Range("E8").Value = "  -0.98%  "
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "'23.08"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "'0.07671"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "1.827.45"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "'4.999"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").Value = "'0.6725"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "'82.86"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "'0.000009154"
$ws.Range("E16").Value = "  -4.80%  "
$ws.Range("D17").Value = "'5.907"
$ws.Range("E17").Value = "  -2.74%  "
$ws.Range("D18").Value = "29.019.40"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "2.079.71"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "'239.75"
$ws.Range("E20").Value = "  +5.86%  "
$ws.Range("D21").Value = "'12.68"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'7.202"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").Value = "'0.9999"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'159.25"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("D26").Value = "'0.1418"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").Value = "'8.489"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").Value = "'17.86"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "'1.498"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "'4.146"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "'4.119"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("D33").Value = "'1.204"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").Value = "'0.7395"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Value = "'2.657"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "'2.775"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").Value = "1.208.85"
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("D41").Value = "'6.449"
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("D42").Value = "'0.8978"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("D43").Value = "'0.9996"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "'101.60"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "1.977.84"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'65.54"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").Value = "'0.5084"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").Value = "  -4.25%  "
$ws.Range("D49").Value = "'0.4068"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "'9.114"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("E51").Value = "  +0.56%  "
